$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 31086.5
$ws.Range("J87").Value = 31086.5
$ws.Range("L87").Value = 31086.5
$ws.Range("N87").Value = -33582.5
$ws.Range("H90").Value = 31086.5
$ws.Range("J90").Value = 31086.5
$ws.Range("L90").Value = 93259.5
$ws.Range("N90").Value = -105739.5
$ws.Range("H123").Value = 28000
$ws.Range("J123").Value = 28000
$ws.Range("L123").Value = 28000
$ws.Range("N123").Value = -37800
$ws.Range("H124").Value = 48511.668
$ws.Range("J124").Value = 48511.668
$ws.Range("L124").Value = 48511.668
$ws.Range("N124").Value = -58331.668
$ws.Range("H126").Value = 42996.668
$ws.Range("J126").Value = 42996.668
$ws.Range("L126").Value = 42996.668
$ws.Range("N126").Value = -52876.668
$ws.Range("H128").Value = 34801.2
$ws.Range("J128").Value = 34801.2
$ws.Range("L128").Value = 34801.2
$ws.Range("N128").Value = -44761.2

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 54998
$ws.Range("J80").Value = 54998
$ws.Range("L80").Value = 54998
$ws.Range("N80").Value = -56994
$ws.Range("H83").Value = 54998
$ws.Range("J83").Value = 54998
$ws.Range("L83").Value = 164994
$ws.Range("N83").Value = -174978
$ws.Range("H119").Value = 52694
$ws.Range("J119").Value = 52694
$ws.Range("L119").Value = 52694
$ws.Range("N119").Value = -62370
$ws.Range("H121").Value = 37968.332
$ws.Range("J121").Value = 37968.332
$ws.Range("L121").Value = 37968.332
$ws.Range("N121").Value = -41462.332
$ws.Range("H128").Value = 47425
$ws.Range("J128").Value = 47425
$ws.Range("L128").Value = 47425
$ws.Range("N128").Value = -57385
$ws.Range("H131").Value = 51687
$ws.Range("J131").Value = 51687
$ws.Range("L131").Value = 51687
$ws.Range("N131").Value = -61767
$ws.Range("H137").Value = 52060
$ws.Range("J137").Value = 52060
$ws.Range("L137").Value = 52060
$ws.Range("N137").Value = -62260

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 47676
$ws.Range("J108").Value = 47676
$ws.Range("L108").Value = 47676
$ws.Range("N108").Value = -55356
$ws.Range("H122").Value = 34777.332
$ws.Range("J122").Value = 34777.332
$ws.Range("L122").Value = 34777.332
$ws.Range("N122").Value = -44577.332
$ws.Range("H130").Value = 40555.145
$ws.Range("J130").Value = 40555.145
$ws.Range("L130").Value = 40555.145
$ws.Range("N130").Value = -50595.145
$ws.Range("H139").Value = 58547.25
$ws.Range("J139").Value = 58547.25
$ws.Range("L139").Value = 58547.25
$ws.Range("N139").Value = -68827.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 75
$ws.Range("I13").Value = 75
$ws.Range("K13").Value = 75
$ws.Range("M13").Value = 64
$ws.Range("H20").Value = 47462.8
$ws.Range("J20").Value = 47462.8
$ws.Range("L20").Value = 47462.8
$ws.Range("N20").Value = -47934.8
$ws.Range("H21").Value = 1900
$ws.Range("I21").Value = 3000
$ws.Range("J21").Value = 1350
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 1350
$ws.Range("M21").Value = -2765
$ws.Range("N21").Value = -1820
$ws.Range("H30").Value = 47462.8
$ws.Range("J30").Value = 47462.8
$ws.Range("L30").Value = 47462.8
$ws.Range("N30").Value = -47644.8
$ws.Range("H31").Value = 181671.83
$ws.Range("I31").Value = 2030.6842
$ws.Range("J31").Value = 238558.2
$ws.Range("K31").Value = 2030.6842
$ws.Range("L31").Value = 238558.2
$ws.Range("M31").Value = -1735.6842
$ws.Range("N31").Value = -239148.2
$ws.Range("H34").Value = 181671.83
$ws.Range("I34").Value = 2030.6842
$ws.Range("J34").Value = 238558.2
$ws.Range("K34").Value = 2030.6842
$ws.Range("L34").Value = 238558.2
$ws.Range("M34").Value = -1828.6842
$ws.Range("N34").Value = -238962.2
$ws.Range("H81").Value = 38622.375
$ws.Range("J81").Value = 38622.375
$ws.Range("L81").Value = 38622.375
$ws.Range("N81").Value = -40618.375
$ws.Range("H82").Value = 36664
$ws.Range("J82").Value = 36664
$ws.Range("L82").Value = 36664
$ws.Range("N82").Value = -37386
$ws.Range("H84").Value = 38622.375
$ws.Range("J84").Value = 38622.375
$ws.Range("L84").Value = 115867.125
$ws.Range("N84").Value = -125851.125
$ws.Range("H85").Value = 36664
$ws.Range("J85").Value = 36664
$ws.Range("L85").Value = 36664
$ws.Range("N85").Value = -39160
$ws.Range("H88").Value = 28781.75
$ws.Range("J88").Value = 28781.75
$ws.Range("L88").Value = 28781.75
$ws.Range("N88").Value = -29593.75
$ws.Range("H91").Value = 28781.75
$ws.Range("J91").Value = 28781.75
$ws.Range("L91").Value = 28781.75
$ws.Range("N91").Value = -31589.75
$ws.Range("H128").Value = 47462.8
$ws.Range("J128").Value = 47462.8
$ws.Range("L128").Value = 47462.8
$ws.Range("N128").Value = -57422.8
$ws.Range("H133").Value = 11007
$ws.Range("J133").Value = 11007
$ws.Range("L133").Value = 11007
$ws.Range("N133").Value = -16067

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 46462.1
$ws.Range("J130").Value = 46462.1
$ws.Range("L130").Value = 46462.1
$ws.Range("N130").Value = -56502.1
$ws.Range("H137").Value = 35558
$ws.Range("J137").Value = 35558
$ws.Range("L137").Value = 35558
$ws.Range("N137").Value = -45758

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 38137
$ws.Range("J92").Value = 38137
$ws.Range("L92").Value = 38137
$ws.Range("N92").Value = -43129
$ws.Range("H96").Value = 29346.5
$ws.Range("J96").Value = 29346.5
$ws.Range("L96").Value = 29346.5
$ws.Range("N96").Value = -34838.5
$ws.Range("H111").Value = 39916.5
$ws.Range("J111").Value = 39916.5
$ws.Range("L111").Value = 39916.5
$ws.Range("N111").Value = -48096.5
$ws.Range("H121").Value = 24801
$ws.Range("J121").Value = 24801
$ws.Range("L121").Value = 24801
$ws.Range("N121").Value = -28295
$ws.Range("H124").Value = 42996
$ws.Range("J124").Value = 42996
$ws.Range("L124").Value = 42996
$ws.Range("N124").Value = -52816
$ws.Range("H125").Value = 49715
$ws.Range("J125").Value = 49715
$ws.Range("L125").Value = 49715
$ws.Range("N125").Value = -59555
$ws.Range("H127").Value = 41161.5
$ws.Range("J127").Value = 41161.5
$ws.Range("L127").Value = 41161.5
$ws.Range("N127").Value = -51081.5
$ws.Range("H128").Value = 35206.5
$ws.Range("J128").Value = 35206.5
$ws.Range("L128").Value = 35206.5
$ws.Range("N128").Value = -45166.5
$ws.Range("H139").Value = 84924.5
$ws.Range("J139").Value = 59899.332
$ws.Range("L139").Value = 59899.332
$ws.Range("N139").Value = -70179.33199999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 29990
$ws.Range("J87").Value = 29990
$ws.Range("L87").Value = 29990
$ws.Range("N87").Value = -32486
$ws.Range("H90").Value = 29990
$ws.Range("J90").Value = 29990
$ws.Range("L90").Value = 89970
$ws.Range("N90").Value = -102450
$ws.Range("H92").Value = 29650
$ws.Range("J92").Value = 29650
$ws.Range("L92").Value = 29650
$ws.Range("N92").Value = -34642
$ws.Range("H93").Value = 36397.145
$ws.Range("J93").Value = 36397.145
$ws.Range("L93").Value = 36397.145
$ws.Range("N93").Value = -41389.145
$ws.Range("H99").Value = 37825.09
$ws.Range("J99").Value = 39222.5
$ws.Range("L99").Value = 39222.5
$ws.Range("N99").Value = -45212.5
$ws.Range("H106").Value = 32818.668
$ws.Range("J106").Value = 32818.668
$ws.Range("L106").Value = 32818.668
$ws.Range("N106").Value = -35342.668
$ws.Range("H118").Value = 29895.8
$ws.Range("J118").Value = 33619.75
$ws.Range("L118").Value = 33619.75
$ws.Range("N118").Value = -36933.75
$ws.Range("H119").Value = 45830.332
$ws.Range("J119").Value = 45830.332
$ws.Range("L119").Value = 45830.332
$ws.Range("N119").Value = -55506.332
$ws.Range("H123").Value = 43421
$ws.Range("J123").Value = 43421
$ws.Range("L123").Value = 43421
$ws.Range("N123").Value = -53221
$ws.Range("H125").Value = 39707
$ws.Range("J125").Value = 39707
$ws.Range("L125").Value = 39707
$ws.Range("N125").Value = -49547
$ws.Range("H139").Value = 40618.332
$ws.Range("J139").Value = 40618.332
$ws.Range("L139").Value = 40618.332
$ws.Range("N139").Value = -50898.332
